$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 45600
$ws.Range("A4").Value = 45600
$ws.Range("A5").Value = 45600

$ws.Range("A7").Value = 45604
$ws.Range("A8").Value = 45604
$ws.Range("A9").Value = 45604

$ws.Range("A10").Value = 45605
$ws.Range("A11").Value = 45605
$ws.Range("A12").Value = 45605
$ws.Range("A13").Value = 45605

$ws.Range("A13").Select()
